# Actualizacion automatica del mapa (2025-11-12 12:51:11)
#
# - Fills in the previously-blank "OT" (column E) values for the two
#   existing rows (86 and 88) that were pending an OT number.
# - Appends two brand-new claim rows (89 and 90) at the bottom of the
#   "Optical_Power" sheet, extending the used range from A1:R88 to A1:R90.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (so digit-only / date-looking
# strings like "810712796" or "11/11/2025" are not silently reinterpreted
# by Excel as a Number/Date) while still leaving the cell on the sheet's
# default "Normal" style, exactly like every other text cell in this
# sheet.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Existing rows: fill in the OT number that was previously blank.
# ---------------------------------------------------------------------
Set-TextCell 86 5 "810712796"
Set-TextCell 88 5 "810712875"

# ---------------------------------------------------------------------
# New row 89
# ---------------------------------------------------------------------
Set-TextCell 89 1 "7746"
Set-TextCell 89 2 "11/11/2025"
$ws.Cells.Item(89, 3).Value = "CAMPICHUELO 229"
Set-TextCell 89 4 "6"
Set-TextCell 89 5 "810712887"
$ws.Cells.Item(89, 6).Value = "Optical Power"
$ws.Cells.Item(89, 7).Value = "Pendiente"
$ws.Cells.Item(89, 8).Value = "Cambiar"
$ws.Cells.Item(89, 9).Value = 1
$ws.Cells.Item(89, 10).Value = "Cambio"
$ws.Cells.Item(89, 11).Value = "Sin equipos"
$ws.Cells.Item(89, 12).Value = "Pasante"
$ws.Cells.Item(89, 13).Value = -58.433855
$ws.Cells.Item(89, 14).Value = -34.614487
$ws.Cells.Item(89, 15).Value = "Almagro"
$ws.Cells.Item(89, 16).Value = "Capital Sur"
$ws.Cells.Item(89, 17).Value = "ALM-J"
$ws.Cells.Item(89, 18).Value = "Fuera de Poligono OVL"

# ---------------------------------------------------------------------
# New row 90
# ---------------------------------------------------------------------
Set-TextCell 90 1 "7842"
Set-TextCell 90 2 "11/12/2025"
$ws.Cells.Item(90, 3).Value = "FERRARI 410"
Set-TextCell 90 4 "15"
Set-TextCell 90 5 "810713039"
$ws.Cells.Item(90, 6).Value = "Optical Power"
$ws.Cells.Item(90, 7).Value = "Pendiente"
$ws.Cells.Item(90, 8).Value = "Picada"
$ws.Cells.Item(90, 9).Value = 1
$ws.Cells.Item(90, 10).Value = "Cambio"
$ws.Cells.Item(90, 11).Value = "Sin equipos"
$ws.Cells.Item(90, 12).Value = "Pasante"
$ws.Cells.Item(90, 13).Value = -58.441198
$ws.Cells.Item(90, 14).Value = -34.605341
$ws.Cells.Item(90, 15).Value = "Paternal"
$ws.Cells.Item(90, 16).Value = "Capital Norte"
$ws.Cells.Item(90, 17).Value = "ALM-O"
$ws.Cells.Item(90, 18).Value = "Fuera de Poligono OVL"

Write-Output "Applied map update: filled E86/E88 OT numbers and appended rows 89-90 (dimension now A1:R90)."
